$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-13 with refreshed TPM-derived values
# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Epo"
$ws.Range("C2").Value = "Ephb4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.12382
$ws.Range("H2").Value = 0.37146
$ws.Range("I2").Value = 0.3718173164410999
$ws.Range("J2").Value = 0.3770600965540203
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 26.07194833333334
$ws.Range("N2").Value = 78.215845
$ws.Range("O2").Value = 0.7595928012803946
$ws.Range("P2").Value = 0.7818221335728009
$ws.Range("Q2").Value = 3.228228642633334
$ws.Range("R2").Value = 29.0540577837
$ws.Range("S2").Value = 0.282429756960054
$ws.Range("T2").Value = 0.2947939291730305

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Epo"
$ws.Range("C3").Value = "Ephb4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.12382
$ws.Range("H3").Value = 0.37146
$ws.Range("I3").Value = 0.3718173164410999
$ws.Range("J3").Value = 0.3770600965540203
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.323898
$ws.Range("N3").Value = 15.971694
$ws.Range("O3").Value = 0.155109029208254
$ws.Range("P3").Value = 0.1596482641062294
$ws.Range("Q3").Value = 0.65920505036
$ws.Range("R3").Value = 5.93284545324
$ws.Range("S3").Value = 0.0576722229959972
$ws.Range("T3").Value = 0.0601969898785766

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Epo"
$ws.Range("C4").Value = "Ephb4"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.12382
$ws.Range("H4").Value = 0.37146
$ws.Range("I4").Value = 0.3718173164410999
$ws.Range("J4").Value = 0.3770600965540203
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.927739
$ws.Range("N4").Value = 5.855478
$ws.Range("O4").Value = 0.08529816951135136
$ws.Range("P4").Value = 0.05852960232096958
$ws.Range("Q4").Value = 0.36251264298
$ws.Range("R4").Value = 2.17507585788
$ws.Range("S4").Value = 0.03171533648504871
$ws.Range("T4").Value = 0.0220691775024132

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Epo"
$ws.Range("C5").Value = "Ephb4"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.062059
$ws.Range("H5").Value = 0.186177
$ws.Range("I5").Value = 0.1863560882007609
$ws.Range("J5").Value = 0.1889837872076074
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 26.07194833333334
$ws.Range("N5").Value = 78.215845
$ws.Range("O5").Value = 0.7595928012803946
$ws.Range("P5").Value = 0.7818221335728009
$ws.Range("Q5").Value = 1.617999041618333
$ws.Range("R5").Value = 14.561991374565
$ws.Range("S5").Value = 0.1415547430720723
$ws.Range("T5").Value = 0.1477517077253198

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Epo"
$ws.Range("C6").Value = "Ephb4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.062059
$ws.Range("H6").Value = 0.186177
$ws.Range("I6").Value = 0.1863560882007609
$ws.Range("J6").Value = 0.1889837872076074
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 5.323898
$ws.Range("N6").Value = 15.971694
$ws.Range("O6").Value = 0.155109029208254
$ws.Range("P6").Value = 0.1596482641062294
$ws.Range("Q6").Value = 0.330395785982
$ws.Range("R6").Value = 2.973562073838
$ws.Range("S6").Value = 0.02890551192786779
$ws.Range("T6").Value = 0.03017093357191556

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Epo"
$ws.Range("C7").Value = "Ephb4"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.062059
$ws.Range("H7").Value = 0.186177
$ws.Range("I7").Value = 0.1863560882007609
$ws.Range("J7").Value = 0.1889837872076074
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.927739
$ws.Range("N7").Value = 5.855478
$ws.Range("O7").Value = 0.08529816951135136
$ws.Range("P7").Value = 0.05852960232096958
$ws.Range("Q7").Value = 0.181692554601
$ws.Range("R7").Value = 1.090155327606
$ws.Range("S7").Value = 0.01589583320082085
$ws.Range("T7").Value = 0.011061145910372

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Epo"
$ws.Range("C8").Value = "Ephb4"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.1205036666666667
$ws.Range("H8").Value = 0.361511
$ws.Range("I8").Value = 0.3618587462551512
$ws.Range("J8").Value = 0.366961106351533
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 26.07194833333334
$ws.Range("N8").Value = 78.215845
$ws.Range("O8").Value = 0.7595928012803946
$ws.Range("P8").Value = 0.7818221335728009
$ws.Range("Q8").Value = 3.141765371310556
$ws.Range("R8").Value = 28.275888341795
$ws.Range("S8").Value = 0.2748652987357618
$ws.Range("T8").Value = 0.2868983151059911

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Epo"
$ws.Range("C9").Value = "Ephb4"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.1205036666666667
$ws.Range("H9").Value = 0.361511
$ws.Range("I9").Value = 0.3618587462551512
$ws.Range("J9").Value = 0.366961106351533
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 5.323898
$ws.Range("N9").Value = 15.971694
$ws.Range("O9").Value = 0.155109029208254
$ws.Range("P9").Value = 0.1596482641062294
$ws.Range("Q9").Value = 0.6415492299593334
$ws.Range("R9").Value = 5.773943069634
$ws.Range("S9").Value = 0.05612755884215244
$ws.Range("T9").Value = 0.05858470362352368

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Epo"
$ws.Range("C10").Value = "Ephb4"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.1205036666666667
$ws.Range("H10").Value = 0.361511
$ws.Range("I10").Value = 0.3618587462551512
$ws.Range("J10").Value = 0.366961106351533
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.927739
$ws.Range("N10").Value = 5.855478
$ws.Range("O10").Value = 0.08529816951135136
$ws.Range("P10").Value = 0.05852960232096958
$ws.Range("Q10").Value = 0.352803284543
$ws.Range("R10").Value = 2.116819707258
$ws.Range("S10").Value = 0.03086588867723697
$ws.Range("T10").Value = 0.02147808762201825

# Row 11
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "Epo"
$ws.Range("C11").Value = "Ephb4"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.013891
$ws.Range("H11").Value = 0.027782
$ws.Range("I11").Value = 0.04171308627591115
$ws.Range("J11").Value = 0.02820083885872986
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 26.07194833333334
$ws.Range("N11").Value = 78.215845
$ws.Range("O11").Value = 0.7595928012803946
$ws.Range("P11").Value = 0.7818221335728009
$ws.Range("Q11").Value = 0.3621654342983334
$ws.Range("R11").Value = 2.17299260579
$ws.Range("S11").Value = 0.03168496005437014
$ws.Range("T11").Value = 0.02204804000507493

# Row 12
$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("B12").Value = "Epo"
$ws.Range("C12").Value = "Ephb4"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.013891
$ws.Range("H12").Value = 0.027782
$ws.Range("I12").Value = 0.04171308627591115
$ws.Range("J12").Value = 0.02820083885872986
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 5.323898
$ws.Range("N12").Value = 15.971694
$ws.Range("O12").Value = 0.155109029208254
$ws.Range("P12").Value = 0.1596482641062294
$ws.Range("Q12").Value = 0.073954267118
$ws.Range("R12").Value = 0.443725602708
$ws.Range("S12").Value = 0.006470076317536723
$ws.Range("T12").Value = 0.004502214970135722

# Row 13
$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("B13").Value = "Epo"
$ws.Range("C13").Value = "Ephb4"
$ws.Range("D13").Value = "MuSCs"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.013891
$ws.Range("H13").Value = 0.027782
$ws.Range("I13").Value = 0.04171308627591115
$ws.Range("J13").Value = 0.02820083885872986
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 2.927739
$ws.Range("N13").Value = 5.855478
$ws.Range("O13").Value = 0.08529816951135136
$ws.Range("P13").Value = 0.05852960232096958
$ws.Range("Q13").Value = 0.040669222449
$ws.Range("R13").Value = 0.162676889796
$ws.Range("S13").Value = 0.003558049904004293
$ws.Range("T13").Value = 0.001650583883519205

# Append new rows 14-16 for additional Resolving-Mac sending-cluster records
# Row 14
$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("B14").Value = "Epo"
$ws.Range("C14").Value = "Ephb4"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.01273933333333333
$ws.Range("H14").Value = 0.038218
$ws.Range("I14").Value = 0.03825476282707682
$ws.Range("J14").Value = 0.03879417102810949
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 26.07194833333334
$ws.Range("N14").Value = 78.215845
$ws.Range("O14").Value = 0.7595928012803946
$ws.Range("P14").Value = 0.7818221335728009
$ws.Range("Q14").Value = 0.3321392404677778
$ws.Range("R14").Value = 2.98925316421
$ws.Range("S14").Value = 0.02905804245813639
$ws.Range("T14").Value = 0.0303301415633847

# Row 15
$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("B15").Value = "Epo"
$ws.Range("C15").Value = "Ephb4"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.01273933333333333
$ws.Range("H15").Value = 0.038218
$ws.Range("I15").Value = 0.03825476282707682
$ws.Range("J15").Value = 0.03879417102810949
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 5.323898
$ws.Range("N15").Value = 15.971694
$ws.Range("O15").Value = 0.155109029208254
$ws.Range("P15").Value = 0.1596482641062294
$ws.Range("Q15").Value = 0.06782291125466666
$ws.Range("R15").Value = 0.610406201292
$ws.Range("S15").Value = 0.005933659124699889
$ws.Range("T15").Value = 0.006193422062077856

# Row 16
$ws.Range("A16").Value = "Resolving-Mac"
$ws.Range("B16").Value = "Epo"
$ws.Range("C16").Value = "Ephb4"
$ws.Range("D16").Value = "MuSCs"
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.01273933333333333
$ws.Range("H16").Value = 0.038218
$ws.Range("I16").Value = 0.03825476282707682
$ws.Range("J16").Value = 0.03879417102810949
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 2.927739
$ws.Range("N16").Value = 5.855478
$ws.Range("O16").Value = 0.08529816951135136
$ws.Range("P16").Value = 0.05852960232096958
$ws.Range("Q16").Value = 0.037297443034
$ws.Range("R16").Value = 0.223784658204
$ws.Range("S16").Value = 0.003263061244240541
$ws.Range("T16").Value = 0.002270607402646928
